# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: "002" -> "001" (stays text)
$ws.Range("J2").Value = "'001"

# REPORT_DATE: 2020-06-30 -> 2018-12-31 (text)
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# PARENT_NETPROFIT
$ws.Range("O2").Value = 94163686.68000001
# TOTAL_OPERATE_INCOME
$ws.Range("P2").Value = 890611692.1900001
# TOTAL_OPERATE_COST
$ws.Range("Q2").Value = 771036100.39
# TOE_RATIO (newly populated, was blank)
$ws.Range("R2").Value = 46.0330740519
# OPERATE_COST
$ws.Range("S2").Value = 168240269.76
# OPERATE_EXPENSE
$ws.Range("T2").Value = 168240269.76
# OPERATE_EXPENSE_RATIO (newly populated, was blank)
$ws.Range("U2").Value = 30.2953699038
# SALE_EXPENSE
$ws.Range("V2").Value = 468079784.5
# MANAGE_EXPENSE
$ws.Range("W2").Value = 73333863.04000001
# FINANCE_EXPENSE
$ws.Range("X2").Value = 11929829.04
# OPERATE_PROFIT
$ws.Range("Y2").Value = 121270238.31
# TOTAL_PROFIT
$ws.Range("Z2").Value = 121268550.7
# INCOME_TAX
$ws.Range("AA2").Value = 27104864.02
# OPERATE_TAX_ADD
$ws.Range("AG2").Value = 6775507.39
# TOI_RATIO
$ws.Range("AP2").Value = 44.6173513121
# OPERATE_PROFIT_RATIO
$ws.Range("AQ2").Value = 35.777713408945
# PARENT_NETPROFIT_RATIO
$ws.Range("AR2").Value = 58.547012783481
# DEDUCT_PARENT_NETPROFIT
$ws.Range("AS2").Value = 108296690.67
# DPN_RATIO
$ws.Range("AT2").Value = 47.086314352647
